$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new value is a plain numeric-looking string.
# Excel auto-converts such strings to numbers unless the cell is
# pre-formatted as Text, so mark these as Text first to preserve them
# as literal strings (matching the original "Price" column formatting).
$textCells = @("D4","D5","D6","D7","D14","D19","D20","D22","D24","D28","D29","D30","D35","D36","D39","D40","D41","D44","D46","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values per row.
$ws.Range("D2").Value = "63.228.32"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.663.97"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "608.57"
$ws.Range("E5").Value = "  +4.32%  "
$ws.Range("D6").Value = "143.33"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.663.38"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "27.33"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "63.118.91"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.643.23"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").Value = "11.45"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").Value = "339.80"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "6.85"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "67.60"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "8.54"
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").Value = "541.41"
$ws.Range("E29").Value = "  +17.48%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  +4.97%  "
$ws.Range("E33").Value = "  +7.89%  "
$ws.Range("D34").Value = "0.0₃0809"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").Value = "172.19"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +12.63%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").Value = "19.22"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("D41").Value = "173.73"
$ws.Range("E41").Value = "  +9.20%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "22.03"
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("D46").Value = "0.633"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "18.81"
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("E51").Value = "  -0.70%  "
